$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style of an existing header cell (H1) into I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for I column (all 1 except row 8 = 3)
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 3

# Data values for J column (same as H except row 8 = 4)
$ws.Range("J2").Value = 3
$ws.Range("J3").Value = 5
$ws.Range("J4").Value = 5
$ws.Range("J5").Value = 6
$ws.Range("J6").Value = 3
$ws.Range("J7").Value = 3
$ws.Range("J8").Value = 4
